{"js": "// Revise the SOC111 fitness-center paper: opening hook sentence rework,\n// small wording fixes, and a couple of phrase swaps later in the essay.\n\nconst replacements = [\n  {\n    find: \"The most immediately obvious gender dynamic in the fitness center is the split between the upper and lower floors.\",\n    replace: \"As soon as you walk into the fitness center, the most immediately obvious gender dynamic is that members of each gender tend to keep to one floor.\"\n  },\n  {\n    find: \"as well of most of the center\\u2019s male population\",\n    replace: \"as well as most of the center\\u2019s male population\"\n  },\n  {\n    find: \"My hypothesis on this dynamic as a long time gym-goer myself, is as follows:\",\n    replace: \"My hypothesis on this dynamic, as a long time gym-goer myself, is as follows:\"\n  },\n  {\n    find: \"anecdotal experiences talking to other guys in and around the gym. Meanwhile\",\n    replace: \"anecdotal experiences talking to other guys in the gym. Meanwhile\"\n  },\n  {\n    find: \"I also noticed a distinct different in body language, especially among girls, depending on where in the fitness center they were working out. The girls that I saw in the fitness center were \",\n    replace: \"I also noticed a distinct difference in body language among girls depending on where in the fitness center they were working out. The girls I saw in the fitness center were the most open \"\n  },\n  {\n    find: \"the most open in terms of body language \\u2013 heads up, making eye contact, talking, not looking at phones \\u2013 in the padded/ab workout area of the upper floor. That body language differed from girls\\u2019 body language in the rest of the gym\",\n    replace: \"in terms of body language \\u2013 heads up, making eye contact, talking, not looking at phones \\u2013 in the padded/ab workout area of the upper floor. That was noticeably different from girls\\u2019 body language in the rest of the gym\"\n  },\n  {\n    find: \"working out with a partner. As far as I could tell\",\n    replace: \"working out with others. As far as I could tell\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Revise the SOC111 fitness-center paper: opening hook sentence rework,\n# small wording fixes, and a couple of phrase swaps later in the essay.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{\n        Find    = \"The most immediately obvious gender dynamic in the fitness center is the split between the upper and lower floors.\"\n        Replace = \"As soon as you walk into the fitness center, the most immediately obvious gender dynamic is that members of each gender tend to keep to one floor.\"\n    },\n    @{\n        Find    = \"as well of most of the center\u2019s male population\"\n        Replace = \"as well as most of the center\u2019s male population\"\n    },\n    @{\n        Find    = \"My hypothesis on this dynamic as a long time gym-goer myself, is as follows:\"\n        Replace = \"My hypothesis on this dynamic, as a long time gym-goer myself, is as follows:\"\n    },\n    @{\n        Find    = \"anecdotal experiences talking to other guys in and around the gym. Meanwhile\"\n        Replace = \"anecdotal experiences talking to other guys in the gym. Meanwhile\"\n    },\n    @{\n        Find    = \"I also noticed a distinct different in body language, especially among girls, depending on where in the fitness center they were working out. The girls that I saw in the fitness center were \"\n        Replace = \"I also noticed a distinct difference in body language among girls depending on where in the fitness center they were working out. The girls I saw in the fitness center were the most open \"\n    },\n    @{\n        Find    = \"the most open in terms of body language \u2013 heads up, making eye contact, talking, not looking at phones \u2013 in the padded/ab workout area of the upper floor. That body language differed from girls\u2019 body language in the rest of the gym\"\n        Replace = \"in terms of body language \u2013 heads up, making eye contact, talking, not looking at phones \u2013 in the padded/ab workout area of the upper floor. That was noticeably different from girls\u2019 body language in the rest of the gym\"\n    },\n    @{\n        Find    = \"working out with a partner. As far as I could tell\"\n        Replace = \"working out with others. As far as I could tell\"\n    }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $item.Find\n    $find.Replacement.Text = $item.Replace\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Could not find text: $($item.Find)\"\n    }\n}\n"}
